$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to stay
# text (matching the original inlineStr type) by temporarily switching the
# cell to a Text number format before assigning, then restoring a plain style
# so no stray formatting is left behind.

$ws.Range('D2').Value = '71.145.69'
$ws.Range('E2').Value = '  -2.27%  '
$ws.Range('D3').Value = '3.878.60'
$ws.Range('E3').Value = '  -2.26%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.21'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '168.17'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.59%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.674'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.757'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +4.84%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.65'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.52%  '
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +6.87%  '
$ws.Range('D14').Value = '4.506.46'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.19'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.13%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.869.09'
$ws.Range('E16').Value = '  -2.65%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.84'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('E18').Value = '  -4.27%  '
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').Value = '71.055.95'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '438.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.73'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '94.75'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('E24').Value = '  -4.14%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.96'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('E26').Value = '  -5.74%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.35'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.43'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.27'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.33'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +6.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '13.68'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '48.71'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0996'
$ws.Range('E35').Value = '  +14.83%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '69.75'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '634.30'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.432'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.147'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.29'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.30'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +28.34%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('E44').Value = '  -2.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.16'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.46%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.74'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('E47').Value = '  -3.25%  '
$ws.Range('E48').Value = '  -15.27%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.31'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.861.69'
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('E51').Value = '  +1.80%  '
